$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 0.161797665277291
$ws.Cells.Item(2, 4).Value = 0.0003530347199651729
$ws.Cells.Item(2, 5).Value = 0.07945540271330387
$ws.Cells.Item(2, 6).Value = 0.5977386150767501
$ws.Cells.Item(2, 7).Value = 0.002345103280840417
$ws.Cells.Item(2, 9).Value = 0.3451948428057143
$ws.Cells.Item(2, 13).Value = 0.6505005855176478
$ws.Cells.Item(2, 14).Value = 1.231231664624772
$ws.Cells.Item(2, 15).Value = 1.913375912927279

$ws.Cells.Item(3, 2).Value = 0.1434408339321323
$ws.Cells.Item(3, 4).Value = 0.0003268267869325214
$ws.Cells.Item(3, 5).Value = 0.08169543410162916
$ws.Cells.Item(3, 6).Value = 0.5607374907240228
$ws.Cells.Item(3, 7).Value = 0.002349010242533296
$ws.Cells.Item(3, 9).Value = 0.3459606103821358
$ws.Cells.Item(3, 13).Value = 0.5706096327045316
$ws.Cells.Item(3, 14).Value = 1.195113874855991
$ws.Cells.Item(3, 15).Value = 1.798077163701407

$ws.Cells.Item(4, 2).Value = 0.1321466050026316
$ws.Cells.Item(4, 4).Value = 0.0003109512802952352
$ws.Cells.Item(4, 5).Value = 0.08320312098164528
$ws.Cells.Item(4, 6).Value = 0.5384317102994913
$ws.Cells.Item(4, 7).Value = 0.002351537737991335
$ws.Cells.Item(4, 9).Value = 0.3466510786719397
$ws.Cells.Item(4, 13).Value = 0.5215055623576887
$ws.Cells.Item(4, 14).Value = 1.173347869083756
$ws.Cells.Item(4, 15).Value = 1.728644584535488

$ws.Cells.Item(5, 2).Value = 0.1275387438209634
$ws.Cells.Item(5, 4).Value = 0.0003045347628480854
$ws.Cells.Item(5, 5).Value = 0.0838505088013779
$ws.Cells.Item(5, 6).Value = 0.5294450575587177
$ws.Cells.Item(5, 7).Value = 0.002352600158380138
$ws.Cells.Item(5, 9).Value = 0.3469880226902013
$ws.Cells.Item(5, 13).Value = 0.5014828662995399
$ws.Cells.Item(5, 14).Value = 1.16458265601679
$ws.Cells.Item(5, 15).Value = 1.700689982315339

$ws.Cells.Item(6, 2).Value = 0.1267733001839275
$ws.Cells.Item(6, 4).Value = 0.0003034724573929282
$ws.Cells.Item(6, 5).Value = 0.08395999175061952
$ws.Cells.Item(6, 6).Value = 0.5279590384242283
$ws.Cells.Item(6, 7).Value = 0.002352778535126593
$ws.Cells.Item(6, 9).Value = 0.3470473337663336
$ws.Cells.Item(6, 13).Value = 0.4981573755232489
$ws.Cells.Item(6, 14).Value = 1.163133558604699
$ws.Cells.Item(6, 15).Value = 1.696068586248657

$ws.Cells.Item(7, 2).Value = 0.1320844829555057
$ws.Cells.Item(7, 4).Value = 0.0003108645328904203
$ws.Cells.Item(7, 5).Value = 0.08321171865539334
$ws.Cells.Item(7, 6).Value = 0.5383100965434693
$ws.Cells.Item(7, 7).Value = 0.002351551934661234
$ws.Cells.Item(7, 9).Value = 0.3466553975647599
$ws.Cells.Item(7, 13).Value = 0.5212355789490744
$ws.Cells.Item(7, 14).Value = 1.173229232929458
$ws.Cells.Item(7, 15).Value = 1.7282662068219

$ws.Cells.Item(8, 2).Value = 0.1554732662631721
$ws.Cells.Item(8, 4).Value = 0.00034395239284013
$ws.Cells.Item(8, 5).Value = 0.080200141854462
$ws.Cells.Item(8, 6).Value = 0.5848944421145887
$ws.Cells.Item(8, 7).Value = 0.002346423773312893
$ws.Cells.Item(8, 9).Value = 0.3454132662848295
$ws.Cells.Item(8, 13).Value = 0.6229648977768534
$ws.Cells.Item(8, 14).Value = 1.21869387949954
$ws.Cells.Item(8, 15).Value = 1.873336768882609

$ws.Cells.Item(9, 2).Value = 0.2011390620019995
$ws.Cells.Item(9, 4).Value = 0.000410630195688988
$ws.Cells.Item(9, 5).Value = 0.07535646907792248
$ws.Cells.Item(9, 6).Value = 0.6795620543570777
$ws.Cells.Item(9, 7).Value = 0.002337383030676154
$ws.Cells.Item(9, 9).Value = 0.3447178610087889
$ws.Cells.Item(9, 13).Value = 0.8220479152741689
$ws.Cells.Item(9, 14).Value = 1.311051190562978
$ws.Cells.Item(9, 15).Value = 2.168754100734589

$ws.Cells.Item(10, 2).Value = 0.2345490743204266
$ws.Cells.Item(10, 4).Value = 0.0004608302091995142
$ws.Cells.Item(10, 5).Value = 0.07246214937741868
$ws.Cells.Item(10, 6).Value = 0.7511966446152201
$ws.Cells.Item(10, 7).Value = 0.002331353112524718
$ws.Cells.Item(10, 9).Value = 0.3452586487098941
$ws.Cells.Item(10, 13).Value = 0.9680733196341436
$ws.Cells.Item(10, 14).Value = 1.380791518839345
$ws.Cells.Item(10, 15).Value = 2.392671230362339

$ws.Cells.Item(11, 2).Value = 0.2497138050246974
$ws.Cells.Item(11, 4).Value = 0.0004839591290197731
$ws.Cells.Item(11, 5).Value = 0.07129335664920511
$ws.Cells.Item(11, 6).Value = 0.7842508572332179
$ws.Cells.Item(11, 7).Value = 0.002328741455913663
$ws.Cells.Item(11, 9).Value = 0.3457311362551607
$ws.Cells.Item(11, 13).Value = 1.034454512046736
$ws.Cells.Item(11, 14).Value = 1.412914284820459
$ws.Cells.Item(11, 15).Value = 2.496076052977287

$ws.Cells.Item(12, 2).Value = 0.2554510669212675
$ws.Cells.Item(12, 4).Value = 0.0004927619433168218
$ws.Cells.Item(12, 5).Value = 0.07087233052261865
$ws.Cells.Item(12, 6).Value = 0.796835767170947
$ws.Cells.Item(12, 7).Value = 0.002327771272064991
$ws.Cells.Item(12, 9).Value = 0.3459424571215521
$ws.Cells.Item(12, 13).Value = 1.059584560137921
$ws.Cells.Item(12, 14).Value = 1.425134245695233
$ws.Cells.Item(12, 15).Value = 2.535458041671575

$ws.Cells.Item(13, 2).Value = 0.2542156875014427
$ws.Cells.Item(13, 4).Value = 0.0004908640867995473
$ws.Cells.Item(13, 5).Value = 0.07096204172798615
$ws.Cells.Item(13, 6).Value = 0.7941223446054266
$ws.Cells.Item(13, 7).Value = 0.002327979384202053
$ws.Cells.Item(13, 9).Value = 0.3458955070387475
$ws.Cells.Item(13, 13).Value = 1.054172675005759
$ws.Cells.Item(13, 14).Value = 1.422500000379131
$ws.Cells.Item(13, 15).Value = 2.526966381239447

$ws.Cells.Item(14, 2).Value = 0.2501859216732782
$ws.Cells.Item(14, 4).Value = 0.0004846824373130332
$ws.Cells.Item(14, 5).Value = 0.0712582844920604
$ws.Cells.Item(14, 6).Value = 0.7852848579658058
$ws.Cells.Item(14, 7).Value = 0.002328661262386083
$ws.Cells.Item(14, 9).Value = 0.3457478733574106
$ws.Cells.Item(14, 13).Value = 1.036522120047906
$ws.Cells.Item(14, 14).Value = 1.413918518677661
$ws.Cells.Item(14, 15).Value = 2.499311510739005

$ws.Cells.Item(15, 2).Value = 0.2477168710868511
$ws.Cells.Item(15, 4).Value = 0.0004809018521703123
$ws.Cells.Item(15, 5).Value = 0.07144256006627892
$ws.Cells.Item(15, 6).Value = 0.7798805213727462
$ws.Cells.Item(15, 7).Value = 0.002329081376933807
$ws.Cells.Item(15, 9).Value = 0.3456616579701048
$ws.Cells.Item(15, 13).Value = 1.02570971375782
$ws.Cells.Item(15, 14).Value = 1.408669334507096
$ws.Cells.Item(15, 15).Value = 2.482401467472414

$ws.Cells.Item(16, 2).Value = 0.2335573075359321
$ws.Cells.Item(16, 4).Value = 0.0004593247640434228
$ws.Cells.Item(16, 5).Value = 0.0725415342567608
$ws.Cells.Item(16, 6).Value = 0.7490459454154603
$ws.Cells.Item(16, 7).Value = 0.002331526425613574
$ws.Cells.Item(16, 9).Value = 0.3452323135433133
$ws.Cells.Item(16, 13).Value = 0.9637341839949158
$ws.Cells.Item(16, 14).Value = 1.378700100138445
$ws.Cells.Item(16, 15).Value = 2.385944788613642

$ws.Cells.Item(17, 2).Value = 0.2248619316181077
$ws.Cells.Item(17, 4).Value = 0.0004461644794702124
$ws.Cells.Item(17, 5).Value = 0.07325380658924452
$ws.Cells.Item(17, 6).Value = 0.7302501823430276
$ws.Cells.Item(17, 7).Value = 0.002333059963572956
$ws.Cells.Item(17, 9).Value = 0.345026811972744
$ws.Cells.Item(17, 13).Value = 0.9257020191326575
$ws.Cells.Item(17, 14).Value = 1.360415765230982
$ws.Cells.Item(17, 15).Value = 2.327169150809311

$ws.Cells.Item(18, 2).Value = 0.2198574409008529
$ws.Cells.Item(18, 4).Value = 0.0004386224337968869
$ws.Cells.Item(18, 5).Value = 0.07367738427634407
$ws.Cells.Item(18, 6).Value = 0.7194832906219801
$ws.Cells.Item(18, 7).Value = 0.002333954386278538
$ws.Cells.Item(18, 9).Value = 0.3449299402311077
$ws.Cells.Item(18, 13).Value = 0.9038225695064739
$ws.Cells.Item(18, 14).Value = 1.349936600592741
$ws.Cells.Item(18, 15).Value = 2.293508097480242

$ws.Cells.Item(19, 2).Value = 0.2181624808082177
$ws.Cells.Item(19, 4).Value = 0.0004360734652246734
$ws.Cells.Item(19, 5).Value = 0.07382317789651616
$ws.Cells.Item(19, 6).Value = 0.7158453274354883
$ws.Cells.Item(19, 7).Value = 0.002334259350475632
$ws.Cells.Item(19, 9).Value = 0.3449008092444075
$ws.Cells.Item(19, 13).Value = 0.8964138277995914
$ws.Cells.Item(19, 14).Value = 1.346395019233057
$ws.Cells.Item(19, 15).Value = 2.282135885065657

$ws.Cells.Item(20, 2).Value = 0.2257878975478604
$ws.Cells.Item(20, 4).Value = 0.0004475625614759338
$ws.Cells.Item(20, 5).Value = 0.07317654341116331
$ws.Cells.Item(20, 6).Value = 0.7322464706094536
$ws.Cells.Item(20, 7).Value = 0.002332895436163438
$ws.Cells.Item(20, 9).Value = 0.3450464818902041
$ws.Cells.Item(20, 13).Value = 0.9297510614582762
$ws.Cells.Item(20, 14).Value = 1.362358292851098
$ws.Cells.Item(20, 15).Value = 2.333410876050209

$ws.Cells.Item(21, 2).Value = 0.2513697088198228
$ws.Cells.Item(21, 4).Value = 0.0004864969138687059
$ws.Cells.Item(21, 5).Value = 0.07117068291820239
$ws.Cells.Item(21, 6).Value = 0.7878787901668858
$ws.Cells.Item(21, 7).Value = 0.00232846046873423
$ws.Cells.Item(21, 9).Value = 0.3457903589458908
$ws.Cells.Item(21, 13).Value = 1.041706710733422
$ws.Cells.Item(21, 14).Value = 1.416437605466683
$ws.Cells.Item(21, 15).Value = 2.507428292173074

$ws.Cells.Item(22, 2).Value = 0.2680578436896326
$ws.Cells.Item(22, 4).Value = 0.0005122026983759298
$ws.Cells.Item(22, 5).Value = 0.06998561079454646
$ws.Cells.Item(22, 6).Value = 0.8246344425602814
$ws.Cells.Item(22, 7).Value = 0.002325671458970143
$ws.Cells.Item(22, 9).Value = 0.3464653015882462
$ws.Cells.Item(22, 13).Value = 1.114834827643264
$ws.Cells.Item(22, 14).Value = 1.452105973008685
$ws.Cells.Item(22, 15).Value = 2.622470523858851

$ws.Cells.Item(23, 2).Value = 0.2591540630686495
$ws.Cells.Item(23, 4).Value = 0.0004984584559846184
$ws.Cells.Item(23, 5).Value = 0.07060648566968553
$ws.Cells.Item(23, 6).Value = 0.8049806927178196
$ws.Cells.Item(23, 7).Value = 0.002327150019707281
$ws.Cells.Item(23, 9).Value = 0.3460878538873686
$ws.Cells.Item(23, 13).Value = 1.075808884337533
$ws.Cells.Item(23, 14).Value = 1.43303986916851
$ws.Cells.Item(23, 15).Value = 2.560949349800239

$ws.Cells.Item(24, 2).Value = 0.2253692851471669
$ws.Cells.Item(24, 4).Value = 0.0004469304139718133
$ws.Cells.Item(24, 5).Value = 0.07321143025130006
$ws.Cells.Item(24, 6).Value = 0.7313438271316812
$ws.Cells.Item(24, 7).Value = 0.002332969779222665
$ws.Cells.Item(24, 9).Value = 0.3450375228433202
$ws.Cells.Item(24, 13).Value = 0.9279205337895462
$ws.Cells.Item(24, 14).Value = 1.361479974016987
$ws.Cells.Item(24, 15).Value = 2.330588587835734

$ws.Cells.Item(25, 2).Value = 0.1888085894510425
$ws.Cells.Item(25, 4).Value = 0.0003923883563459185
$ws.Cells.Item(25, 5).Value = 0.07655137570499093
$ws.Cells.Item(25, 6).Value = 0.6535907021999066
$ws.Cells.Item(25, 7).Value = 0.002339720772997057
$ws.Cells.Item(25, 9).Value = 0.3447206318636731
$ws.Cells.Item(25, 13).Value = 0.7682335872939916
$ws.Cells.Item(25, 14).Value = 1.285730266583215
$ws.Cells.Item(25, 15).Value = 2.08764384693194
